# Updated cryptos list on Mon Aug 28 08:41:18 UTC 2023 with GitHub Actions
# Applies per-cell price/volume updates to the crypto price table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.085.21'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '1.649.09'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.26'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5187'
$ws.Range('E6').Value = '  -3.13%  '
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2614'
$ws.Range('E8').Value = '  -1.81%  '
$ws.Range('E9').Value = '  -2.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.40'
$ws.Range('E10').Value = '  -1.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07790'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.457'
$ws.Range('D13').Value = '1.651.76'
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('D14').Value = '1.875.84'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5533'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '0.0₅7964'
$ws.Range('E16').Value = '  -3.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.63'
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('D18').Value = '26.073.23'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.619'
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.56'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.05'
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.938'
$ws.Range('E23').Value = '  -1.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.94'
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1204'
$ws.Range('E26').Value = '  -2.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.156'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('E28').Value = '  -1.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.472'
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05603'
$ws.Range('E30').Value = '  -4.08%  '
$ws.Range('E31').Value = '  -1.58%  '
$ws.Range('E32').Value = '  -4.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.385'
$ws.Range('E33').Value = '  +3.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.593'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.801'
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9458'
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.404'
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5629'
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.948'
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01573'
$ws.Range('E40').Value = '  -2.00%  '
$ws.Range('D41').Value = '1.058.26'
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('E42').Value = '  -0.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8373'
$ws.Range('E43').Value = '  -3.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.68'
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('D45').Value = '1.788.30'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.98'
$ws.Range('E47').Value = '  +2.78%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05325'
$ws.Range('E49').Value = '  +3.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4336'
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.922'
$ws.Range('E51').Value = '  -1.35%  '
